# Apply updated statistics values for the new run of trials.
# Sheet1 holds the raw stats; Sheet2 pulls the first row via INDEX() formulas
# which will recalculate automatically once Sheet1 is updated.

$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item("Sheet1")

$sheet1.Range("A1").Value = -48.242649518285319
$sheet1.Range("C1").Value = 25.467418578267914
$sheet1.Range("E1").Value = 56.051484437544296
$sheet1.Range("G1").Value = 0.82077446204103555
$sheet1.Range("H1").Value = 0.15002041082696027
$sheet1.Range("I1").Value = 0.062873566091793481

$excel.CalculateFullRebuild()
